# Update of league bases, rotating some duplicate/mis-ordered match rows
# into their correct row positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data columns (every column except A, which holds the row's own id
# and must stay untouched).
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

function Read-Row($ws, $cols, $r) {
    $data = @{}
    foreach ($c in $cols) {
        $data[$c] = $ws.Range($c + $r).Value2
    }
    return $data
}

function Write-Row($ws, $cols, $r, $data) {
    foreach ($c in $cols) {
        $ws.Range($c + $r).Value = $data[$c]
    }
}

# --- Rows 95, 96, 97: 3-way cyclic rotation ---
# new row95 = old row97, new row96 = old row95, new row97 = old row96
$row95 = Read-Row $ws $cols 95
$row96 = Read-Row $ws $cols 96
$row97 = Read-Row $ws $cols 97

Write-Row $ws $cols 95 $row97
Write-Row $ws $cols 96 $row95
Write-Row $ws $cols 97 $row96

# --- Rows 215, 216: swap ---
$row215 = Read-Row $ws $cols 215
$row216 = Read-Row $ws $cols 216

Write-Row $ws $cols 215 $row216
Write-Row $ws $cols 216 $row215
